# Normalize the "Recorded By" (column G) values on the active worksheet.
# A handful of distinct "Recorded By" strings need their comma-separated
# names reordered. The mapping below was derived from the target diff and
# is applied to every cell in column G whose text exactly matches one of
# the keys (so it naturally only touches the rows that actually changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
